$d = $word.ActiveDocument

# 1) "Choose a new super stat" -> "Open an additional slot for Super Stats"
$d.Content.Find.Execute("Choose a new super stat", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Open an additional slot for Super Stats", 2)

# 2) "Note: You can only get one super stat per hero tier" -> "You still need to buy the Super Stat for 30P"
$d.Content.Find.Execute("Note: You can only get one super stat per hero tier", $true, $false, $false, $false, $false,
                         $true, 1, $false, "You still need to buy the Super Stat for 30P", 2)

Write-Host "Done part 1/2"
